# Update the "dSF" (F) column values for the specified rows to reflect the
# repulled data / mean calculation described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = -2
    12 = -1
    18 = -9
    21 = 4
    25 = 5
    30 = 2
    32 = 4
    33 = -2
    34 = -2
    41 = -4
    42 = -1
    44 = 4
    45 = 1
    47 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
